$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "UK": append row 54 (new wave 12 / week 42 record for panel E)
# ---------------------------------------------------------------------------
$wsUK = $wb.Worksheets.Item("UK")

$wsUK.Cells.Item(54, 1).Value = 3
$wsUK.Cells.Item(54, 2).Value = 0
$wsUK.Cells.Item(54, 3).Value = "uk"
$wsUK.Cells.Item(54, 4).Value = 42
$wsUK.Cells.Item(54, 5).Value = "E"

# Copy the date cell's number formatting (style) from an existing date cell
# so the new date cell matches the rest of the column's style.
$wsUK.Cells.Item(2, 7).Copy()
$wsUK.Cells.Item(54, 7).PasteSpecial(-4122)
$wsUK.Cells.Item(54, 7).Value = 44210

$wsUK.Cells.Item(54, 6).Formula = "=F52+1"
$wsUK.Cells.Item(54, 8).Value = "20-040199_PEW12_Final_DPClean_IntUse"
$wsUK.Cells.Item(54, 9).Formula = '=C54&"_"&"wk"&TEXT(D54,"00")&"_"&YEAR(G54)&TEXT(G54,"MM")&TEXT(G54,"DD")&"_p"&E54&"_wv"&TEXT(F54,"00")&""'
$wsUK.Cells.Item(54, 10).Value = 1

# ---------------------------------------------------------------------------
# Sheet "BE": row 10 had its spss_name placed in the wrong column (I instead
# of H). Move it to H10 and populate I10 with the r_name formula, like the
# other rows in the sheet.
# ---------------------------------------------------------------------------
$wsBE = $wb.Worksheets.Item("BE")

$wsBE.Cells.Item(10, 8).Value = "20_060765_BE2_Wave1_Final_v1_20112020_IntClientUse"
$wsBE.Cells.Item(10, 9).Formula = '=A10&"_"&"wk"&TEXT(D10,"00")&"_"&YEAR(G10)&TEXT(G10,"MM")&TEXT(G10,"DD")&"_p"&E10&"_wv"&TEXT(F10,"00")&""'

# ---------------------------------------------------------------------------
# Sheet "Group1": add wave 2 rows (week 2) for each of the 7 countries.
# ---------------------------------------------------------------------------
$wsG1 = $wb.Worksheets.Item("Group1")

$countries = @("at", "dk", "es", "fr", "it", "pl", "pt")

for ($i = 0; $i -lt $countries.Length; $i++) {
    $r = 9 + $i

    $wsG1.Cells.Item($r, 1).Value = $countries[$i]
    $wsG1.Cells.Item($r, 2).Value = 5
    $wsG1.Cells.Item($r, 3).Value = 0
    $wsG1.Cells.Item($r, 4).Value = 2
    $wsG1.Cells.Item($r, 6).Value = 2

    $wsG1.Cells.Item(2, 7).Copy()
    $wsG1.Cells.Item($r, 7).PasteSpecial(-4122)
    $wsG1.Cells.Item($r, 7).Value = 44214

    $wsG1.Cells.Item($r, 8).Value = "20-030971_G1_Merged_Wave2_Final_v1_18012021_IntClientUse"
    $wsG1.Cells.Item($r, 9).Formula = '=A' + $r + '&"_"&"wk"&TEXT(D' + $r + ',"00")&"_"&YEAR(G' + $r + ')&TEXT(G' + $r + ',"MM")&TEXT(G' + $r + ',"DD")&"_p"&E' + $r + '&"_wv"&TEXT(F' + $r + ',"00")&""'
}

# ---------------------------------------------------------------------------
# Selections, to match the view state saved with the workbook. Group1 must
# be selected last so it remains the active (tab-selected) sheet.
# ---------------------------------------------------------------------------
$wsUK.Range("A54").Select()
$wsBE.Range("I9:I10").Select()
$wsG1.Range("I8:I15").Select()
